# Updated symbol list on Sun Dec 25 06:52:43 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as literal text (e.g. "3.390",
# "0.03969") rather than numbers, so trailing/format-significant zeros
# survive. We re-enter each updated price with a leading apostrophe
# (quote-prefix) so Excel keeps storing it as text instead of coercing
# it to a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'245.11"
$ws.Range("D3").Value  = "'23.07"
$ws.Range("D4").Value  = "'5.409"
$ws.Range("D5").Value  = "'0.06009"
$ws.Range("D6").Value  = "'3.389"
$ws.Range("D7").Value  = "'0.8096"
$ws.Range("D8").Value  = "'0.9255"
$ws.Range("D9").Value  = "'0.1419"
$ws.Range("D10").Value = "'0.07461"
$ws.Range("D11").Value = "'0.03386"
$ws.Range("D12").Value = "'0.03030"
$ws.Range("D13").Value = "'0.09356"
$ws.Range("D14").Value = "'3.949"
$ws.Range("D15").Value = "'0.001598"
$ws.Range("D16").Value = "'0.04835"
$ws.Range("D18").Value = "'0.005270"
$ws.Range("D20").Value = "'0.0009859"
$ws.Range("D22").Value = "'3.653"
$ws.Range("D23").Value = "'6.439"
$ws.Range("D40").Value = "'0.03970"

$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("D42").Value = "'0.1074"
$ws.Range("D43").Value = "'0.002711"
$ws.Range("D44").Value = "'0.006711"
$ws.Range("D45").Value = "'0.00005196"

$ws.Range("E48").Value = "47CoinbaseStockTokenCOINBestin24h"

$ws.Range("D49").Value = "'0.002012"
